$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws.Cells.Item(2, 4) '66.421.07'
Set-TextValue $ws.Cells.Item(2, 5) '  -0.41%  '

Set-TextValue $ws.Cells.Item(3, 4) '3.331.49'
Set-TextValue $ws.Cells.Item(3, 5) '  -0.67%  '

Set-TextValue $ws.Cells.Item(4, 5) '  -0.03%  '

Set-TextValue $ws.Cells.Item(5, 4) '586.72'
Set-TextValue $ws.Cells.Item(5, 5) '  +2.21%  '

Set-TextValue $ws.Cells.Item(6, 4) '182.19'
Set-TextValue $ws.Cells.Item(6, 5) '  +0.08%  '

Set-TextValue $ws.Cells.Item(7, 4) '0.652'
Set-TextValue $ws.Cells.Item(7, 5) '  +4.09%  '

Set-TextValue $ws.Cells.Item(8, 5) '  -0.02%  '

Set-TextValue $ws.Cells.Item(9, 4) '3.331.67'
Set-TextValue $ws.Cells.Item(9, 5) '  -0.66%  '

Set-TextValue $ws.Cells.Item(10, 4) '0.127'
Set-TextValue $ws.Cells.Item(10, 5) '  -1.85%  '

Set-TextValue $ws.Cells.Item(11, 5) '  +2.41%  '

Set-TextValue $ws.Cells.Item(12, 4) '0.405'
Set-TextValue $ws.Cells.Item(12, 5) '  +0.06%  '

Set-TextValue $ws.Cells.Item(13, 4) '3.909.05'
Set-TextValue $ws.Cells.Item(13, 5) '  -0.92%  '

Set-TextValue $ws.Cells.Item(14, 5) '  -2.34%  '

Set-TextValue $ws.Cells.Item(15, 4) '66.428.60'
Set-TextValue $ws.Cells.Item(15, 5) '  -0.58%  '

Set-TextValue $ws.Cells.Item(16, 4) '26.56'
Set-TextValue $ws.Cells.Item(16, 5) '  -1.43%  '

Set-TextValue $ws.Cells.Item(17, 4) '0.0000166'
Set-TextValue $ws.Cells.Item(17, 5) '  -0.95%  '

Set-TextValue $ws.Cells.Item(18, 4) '3.299.62'
Set-TextValue $ws.Cells.Item(18, 5) '  -1.80%  '

Set-TextValue $ws.Cells.Item(19, 4) '425.96'
Set-TextValue $ws.Cells.Item(19, 5) '  -2.52%  '

Set-TextValue $ws.Cells.Item(20, 4) '5.55'
Set-TextValue $ws.Cells.Item(20, 5) '  -2.55%  '

Set-TextValue $ws.Cells.Item(21, 4) '13.17'
Set-TextValue $ws.Cells.Item(21, 5) '  -3.37%  '

Set-TextValue $ws.Cells.Item(22, 4) '7.41'
Set-TextValue $ws.Cells.Item(22, 5) '  -2.63%  '

Set-TextValue $ws.Cells.Item(23, 4) '71.95'
Set-TextValue $ws.Cells.Item(23, 5) '  -2.17%  '

Set-TextValue $ws.Cells.Item(24, 5) '  +0.15%  '

Set-TextValue $ws.Cells.Item(25, 5) '  +0.22%  '

Set-TextValue $ws.Cells.Item(26, 4) '3.464.84'
Set-TextValue $ws.Cells.Item(26, 5) '  -1.11%  '

Set-TextValue $ws.Cells.Item(27, 4) '0.516'
Set-TextValue $ws.Cells.Item(27, 5) '  -0.63%  '

Set-TextValue $ws.Cells.Item(28, 5) '  +4.80%  '

Set-TextValue $ws.Cells.Item(29, 4) '0.0000115'
Set-TextValue $ws.Cells.Item(29, 5) '  -1.33%  '

Set-TextValue $ws.Cells.Item(30, 4) '9.03'
Set-TextValue $ws.Cells.Item(30, 5) '  -0.60%  '

Set-TextValue $ws.Cells.Item(31, 4) '0.999'
Set-TextValue $ws.Cells.Item(31, 5) '  -0.02%  '

Set-TextValue $ws.Cells.Item(32, 5) '  -1.63%  '

Set-TextValue $ws.Cells.Item(33, 4) '22.44'
Set-TextValue $ws.Cells.Item(33, 5) '  -2.03%  '

Set-TextValue $ws.Cells.Item(35, 4) '5.21'
Set-TextValue $ws.Cells.Item(35, 5) '  -1.71%  '

Set-TextValue $ws.Cells.Item(36, 4) '6.63'
Set-TextValue $ws.Cells.Item(36, 5) '  -2.50%  '

Set-TextValue $ws.Cells.Item(37, 4) '1.19'
Set-TextValue $ws.Cells.Item(37, 5) '  -2.99%  '

Set-TextValue $ws.Cells.Item(38, 4) '160.68'
Set-TextValue $ws.Cells.Item(38, 5) '  -0.30%  '

Set-TextValue $ws.Cells.Item(39, 4) '1.44'
Set-TextValue $ws.Cells.Item(39, 5) '  -2.67%  '

Set-TextValue $ws.Cells.Item(40, 2) 'Stacks'
Set-TextValue $ws.Cells.Item(40, 3) 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Cells.Item(40, 4) '1.81'
Set-TextValue $ws.Cells.Item(40, 5) '  +0.78%  '

Set-TextValue $ws.Cells.Item(41, 2) 'Maker'
Set-TextValue $ws.Cells.Item(41, 3) 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Cells.Item(41, 4) '2.873.25'
Set-TextValue $ws.Cells.Item(41, 5) '  +1.93%  '

Set-TextValue $ws.Cells.Item(42, 4) '26.55'
Set-TextValue $ws.Cells.Item(42, 5) '  -5.40%  '

Set-TextValue $ws.Cells.Item(43, 4) '4.34'
Set-TextValue $ws.Cells.Item(43, 5) '  -2.47%  '

Set-TextValue $ws.Cells.Item(44, 4) '0.762'
Set-TextValue $ws.Cells.Item(44, 5) '  -4.92%  '

Set-TextValue $ws.Cells.Item(45, 4) '39.85'
Set-TextValue $ws.Cells.Item(45, 5) '  -1.62%  '

Set-TextValue $ws.Cells.Item(46, 4) '0.0665'
Set-TextValue $ws.Cells.Item(46, 5) '  -0.83%  '

Set-TextValue $ws.Cells.Item(47, 4) '5.95'
Set-TextValue $ws.Cells.Item(47, 5) '  -4.25%  '

Set-TextValue $ws.Cells.Item(48, 4) '2.32'
Set-TextValue $ws.Cells.Item(48, 5) '  -1.21%  '

Set-TextValue $ws.Cells.Item(49, 4) '23.20'
Set-TextValue $ws.Cells.Item(49, 5) '  -4.70%  '

Set-TextValue $ws.Cells.Item(50, 4) '314.35'
Set-TextValue $ws.Cells.Item(50, 5) '  -3.63%  '

Set-TextValue $ws.Cells.Item(51, 4) '0.0273'
Set-TextValue $ws.Cells.Item(51, 5) '  -0.09%  '
